$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - "Data driving": add a note about clunky early components
$ws.Range("D8").Value = "Can be clunky for earlier components built while I was still learning."

# Row 13 - "debug drawing": now implemented ("Yes"), add usage note
$ws.Range("B13").Value = "Yes"
$ws.Range("D13").Value = "Ingame: Pause menu->Debug->Toggle Collision Drawing"

# Row 19 - "sprite z sorting": partially implemented, update status/location/note
$ws.Range("B19").Value = "Sort-of"
$ws.Range("C19").Value = "GraphicsManager"
$ws.Range("D19").Value = "I have four render layers. Not using any depth buffer as that was unnecessary, but this does do some sorting."

# Row 23 - "advanced collision response": now implemented ("Yes"), add note
$ws.Range("B23").Value = "Yes"
$ws.Range("D23").Value = "Most collisions will be 'piercing', though this can be modified in config.json"

# Row 24 - "Two distinct AI behaviors": now implemented ("Yes"), add location
$ws.Range("B24").Value = "Yes"
$ws.Range("C24").Value = "AIEnemyCore, AIEnemyStationary, AIEnemyChase, AIMissile"

# Row 25 - "Two distinct weapons": now implemented ("Yes"), add location + note (note lands in column F)
$ws.Range("B25").Value = "Yes"
$ws.Range("C25").Value = "Turret, AIMissile, ControllerShip, Events.h, FollowCursor, TargetLock"
$ws.Range("F25").Value = "Press 'T' to target an enemy, then press space to launch a missile at them."

# Reflect the post-edit selection state
$ws.Range("G28").Select()
